# Generate Report for Handoff
# Adds a new handed-off file (acf846a4-64d9-4b75-bcd5-5fbcabd4917e) as row 3
# on the Overview, zh-cn and de-de sheets, mirroring the existing
# 68a64f2b-bdf0-426f-99d6-5555a6810d97 row.

$wb = $excel.ActiveWorkbook

$newGuid  = "acf846a4-64d9-4b75-bcd5-5fbcabd4917e"
$mdName   = "$newGuid.md"
$hashZh   = "ae39481b34e0c96000a53fff63269ef3bfac2bc9"
$xlfZh    = "$newGuid.$hashZh.zh-cn.xlf"
$xlfDe    = "$newGuid.$hashZh.de-de.xlf"

$mdUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/c20c2daff9a9ea0e6dab474caad9d9c4f0c8635d/e2e/$mdName"
$zhUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8e9a9557c0fc5331dc8c13f5fc9abd43afabca1e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$xlfZh"
$deUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9d2648d0e0c27ac4c644bd7bdd2c36680b91e1a4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$xlfDe"

$handoffDateTime = "2016-03-23 14:39:07"
$zhHandoffDateTime = "2016-03-23 14:39:00"
$deHandoffDateTime = "2016-03-23 14:39:07"
$epoch = "0001-01-01 00:00:00"

# ---------------------------------------------------------------
# Overview sheet: File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = $mdName
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdName)
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = $handoffDateTime
$wsOverview.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A3").Value = $mdName
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdName)
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = $xlfZh
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), $zhUrl, [Type]::Missing, [Type]::Missing, $xlfZh)
$wsZh.Range("E3").Value = $zhHandoffDateTime
$wsZh.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("H3").Value = $epoch
$wsZh.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("J3").Value = "Include"

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A3").Value = $mdName
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdName)
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = $xlfDe
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), $deUrl, [Type]::Missing, [Type]::Missing, $xlfDe)
$wsDe.Range("E3").Value = $deHandoffDateTime
$wsDe.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("H3").Value = $epoch
$wsDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("J3").Value = "Include"
